# Applies scheduled-runner price/profit refresh to Sheets/Louisoix_Profits.xlsx
# (workbook tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# ALC
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H43").Value = 3203.2693
$ws.Range("I43").Value = 2269.5715
$ws.Range("J43").Value = 3547.2632
$ws.Range("K43").Value = 2269.5715
$ws.Range("L43").Value = 3547.2632
$ws.Range("M43").Value = -2200.5715
$ws.Range("N43").Value = -3685.2632

$ws.Range("H74").Value = 21473.5
$ws.Range("I74").Value = 8822
$ws.Range("K74").Value = 8822
$ws.Range("M74").Value = -7886

$ws.Range("H77").Value = 21473.5
$ws.Range("I77").Value = 8822
$ws.Range("K77").Value = 44110
$ws.Range("M77").Value = -39430

$ws.Range("H132").Value = 5520.0684
$ws.Range("I132").Value = 5299.6045
$ws.Range("K132").Value = 15898.8135
$ws.Range("M132").Value = -13368.8135

# ---------------------------------------------------------------------------
# ARM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 34299.574
$ws.Range("I32").Value = 37620.965
$ws.Range("K32").Value = 37620.965
$ws.Range("M32").Value = -37333.965

$ws.Range("H61").Value = 4666.2
$ws.Range("I61").Value = 3956.125
$ws.Range("J61").Value = 7506.5
$ws.Range("K61").Value = 3956.125
$ws.Range("L61").Value = 7506.5
$ws.Range("M61").Value = -3744.125
$ws.Range("N61").Value = -7930.5

$ws.Range("H63").Value = 6378.4443
$ws.Range("I63").Value = 6058.2856
$ws.Range("K63").Value = 6058.2856
$ws.Range("M63").Value = -5372.2856

$ws.Range("H66").Value = 6378.4443
$ws.Range("I66").Value = 6058.2856
$ws.Range("K66").Value = 30291.428
$ws.Range("M66").Value = -26859.428

$ws.Range("H88").Value = 2714.6843
$ws.Range("I88").Value = 1799.1666
$ws.Range("J88").Value = 3137.2307
$ws.Range("K88").Value = 1799.1666
$ws.Range("L88").Value = 3137.2307
$ws.Range("M88").Value = -1393.1666
$ws.Range("N88").Value = -3949.2307

$ws.Range("H91").Value = 2714.6843
$ws.Range("I91").Value = 1799.1666
$ws.Range("J91").Value = 3137.2307
$ws.Range("K91").Value = 1799.1666
$ws.Range("L91").Value = 3137.2307
$ws.Range("M91").Value = -395.1666
$ws.Range("N91").Value = -5945.2307

$ws.Range("H96").Value = 26000
$ws.Range("J96").Value = 26000
$ws.Range("L96").Value = 26000
$ws.Range("N96").Value = -31492

$ws.Range("H97").Value = 687.5833
$ws.Range("J97").Value = 817.1818
$ws.Range("L97").Value = 817.1818
$ws.Range("N97").Value = -1809.1818

$ws.Range("H136").Value = 4666.2
$ws.Range("I136").Value = 3956.125
$ws.Range("J136").Value = 7506.5
$ws.Range("K136").Value = 11868.375
$ws.Range("L136").Value = 22519.5
$ws.Range("M136").Value = -9318.375
$ws.Range("N136").Value = -27619.5

# ---------------------------------------------------------------------------
# BSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H94").Value = 1725.6786
$ws.Range("I94").Value = 1560.125
$ws.Range("K94").Value = 1560.125
$ws.Range("M94").Value = -1109.125

$ws.Range("H105").Value = 3982
$ws.Range("I105").Value = 3876.9285
$ws.Range("J105").Value = 4349.75
$ws.Range("K105").Value = 3876.9285
$ws.Range("L105").Value = 4349.75
$ws.Range("M105").Value = -2129.9285
$ws.Range("N105").Value = -7843.75

# ---------------------------------------------------------------------------
# GSM
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").Value = ""

$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").Value = ""

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = ""

$ws.Range("H64").Value = 53415.332
$ws.Range("I64").Value = 39987.5
$ws.Range("K64").Value = 39987.5
$ws.Range("M64").Value = -39739.5

$ws.Range("H67").Value = 53415.332
$ws.Range("I67").Value = 39987.5
$ws.Range("K67").Value = 39987.5
$ws.Range("M67").Value = -39129.5

$ws.Range("H80").Value = 8740.546
$ws.Range("I80").Value = 2999.75
$ws.Range("J80").Value = 12021
$ws.Range("K80").Value = 2999.75
$ws.Range("L80").Value = 12021
$ws.Range("M80").Value = -2001.75
$ws.Range("N80").Value = -14017

$ws.Range("H83").Value = 8740.546
$ws.Range("I83").Value = 2999.75
$ws.Range("J83").Value = 12021
$ws.Range("K83").Value = 14998.75
$ws.Range("L83").Value = 60105
$ws.Range("M83").Value = -10006.75
$ws.Range("N83").Value = -70089

$ws.Range("H97").Value = 1397.6
$ws.Range("J97").Value = 1236
$ws.Range("L97").Value = 1236
$ws.Range("N97").Value = -2228

$ws.Range("H107").Value = 100663
$ws.Range("I107").Value = 167038.5
$ws.Range("J107").Value = 1099.75
$ws.Range("K107").Value = 167038.5
$ws.Range("L107").Value = 1099.75
$ws.Range("M107").Value = -165118.5
$ws.Range("N107").Value = -4939.75

$ws.Range("H113").Value = 122805.82
$ws.Range("I113").Value = 88023.914
$ws.Range("K113").Value = 88023.914
$ws.Range("M113").Value = -85853.914

# ---------------------------------------------------------------------------
# LTW
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H69").Value = 70000
$ws.Range("J69").Value = 70000
$ws.Range("L69").Value = 70000
$ws.Range("N69").Value = -71622

$ws.Range("H72").Value = 70000
$ws.Range("J72").Value = 70000
$ws.Range("L72").Value = 210000
$ws.Range("N72").Value = -218112

$ws.Range("H92").Value = 89999
$ws.Range("J92").Value = 89999
$ws.Range("L92").Value = 89999
$ws.Range("N92").Value = -94991

$ws.Range("H136").Value = 5818.5
$ws.Range("I136").Value = 1059.6666
$ws.Range("K136").Value = 3178.9998
$ws.Range("M136").Value = -628.9998000000001

$ws.Range("H139").Value = 95499
$ws.Range("J139").Value = 95499
$ws.Range("L139").Value = 95499
$ws.Range("N139").Value = -105779

# ---------------------------------------------------------------------------
# WVR
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H81").Value = 869
$ws.Range("I81").Value = 843.6667
$ws.Range("J81").Value = 945
$ws.Range("K81").Value = 1687.3334
$ws.Range("L81").Value = 1890
$ws.Range("M81").Value = -626.3334
$ws.Range("N81").Value = -4012

$ws.Range("H84").Value = 869
$ws.Range("I84").Value = 843.6667
$ws.Range("J84").Value = 945
$ws.Range("K84").Value = 8436.666999999999
$ws.Range("L84").Value = 9450
$ws.Range("M84").Value = -3132.666999999999
$ws.Range("N84").Value = -20058

$ws.Range("H99").Value = 69999
$ws.Range("J99").Value = 69999
$ws.Range("L99").Value = 69999
$ws.Range("N99").Value = -75989

$ws.Range("H127").Value = 99134.664
$ws.Range("J127").Value = 99134.664
$ws.Range("L127").Value = 99134.664
$ws.Range("N127").Value = -109054.664
